# Update the multiplication answer table cells to match the newly
# generated problem set (commit c8c62b6).
$d = $word.ActiveDocument

$replacements = @(
    @("16×88=1408", "58×39=2262"),
    @("24×19=456",  "50×57=2850"),
    @("15×29=435",  "16×49=784"),
    @("92×12=1104", "14×53=742"),
    @("97×59=5723", "12×60=720"),
    @("93×73=6789", "44×33=1452"),
    @("99×61=6039", "39×81=3159"),
    @("99×69=6831", "84×99=8316"),
    @("20×32=640",  "91×24=2184"),
    @("16×50=800",  "31×84=2604"),
    @("20×88=1760", "37×79=2923"),
    @("69×42=2898", "13×19=247"),
    @("65×15=975",  "57×80=4560"),
    @("26×88=2288", "92×29=2668"),
    @("34×18=612",  "61×63=3843"),
    @("12×17=204",  "68×35=2380"),
    @("79×25=1975", "70×32=2240"),
    @("67×36=2412", "76×97=7372"),
    @("97×99=9603", "88×41=3608"),
    @("47×61=2867", "65×67=4355"),
    @("77×11=847",  "91×25=2275"),
    @("39×92=3588", "23×42=966"),
    @("72×49=3528", "76×88=6688"),
    @("72×54=3888", "56×20=1120"),
    @("56×33=1848", "41×96=3936")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
